$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Apr 08 18:04:51 EDT 2024"
$ws.Range("B3").Value = "Mon Apr 08 18:05:04 EDT 2024"
$ws.Range("B4").Value = "Mon Apr 08 18:05:16 EDT 2024"
$ws.Range("B5").Value = "Mon Apr 08 18:05:29 EDT 2024"
$ws.Range("B6").Value = "Mon Apr 08 18:05:42 EDT 2024"
$ws.Range("B7").Value = "Mon Apr 08 18:05:54 EDT 2024"
